$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '66.718.87'
$ws.Cells.Item(2,4).ClearFormats()
$ws.Cells.Item(2,5).Value = '  +2.34%  '

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '3.696.76'
$ws.Cells.Item(3,4).ClearFormats()
$ws.Cells.Item(3,5).Value = '  +4.83%  '

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '0.998'
$ws.Cells.Item(4,4).ClearFormats()
$ws.Cells.Item(4,5).Value = '  -0.14%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '419.15'
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).Value = '  -0.73%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '129.86'
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).Value = '  -1.22%  '

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '3.689.23'
$ws.Cells.Item(7,4).ClearFormats()
$ws.Cells.Item(7,5).Value = '  +4.72%  '

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.640'
$ws.Cells.Item(8,4).ClearFormats()
$ws.Cells.Item(8,5).Value = '  +0.02%  '

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.999'
$ws.Cells.Item(9,4).ClearFormats()
$ws.Cells.Item(9,5).Value = '  +0.00%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.761'
$ws.Cells.Item(10,4).ClearFormats()
$ws.Cells.Item(10,5).Value = '  -3.60%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.181'
$ws.Cells.Item(11,4).ClearFormats()
$ws.Cells.Item(11,5).Value = '  +7.72%  '

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.0000392'
$ws.Cells.Item(12,4).ClearFormats()
$ws.Cells.Item(12,5).Value = '  +43.68%  '

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '42.89'
$ws.Cells.Item(13,4).ClearFormats()
$ws.Cells.Item(13,5).Value = '  -1.57%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '10.60'
$ws.Cells.Item(14,4).ClearFormats()
$ws.Cells.Item(14,5).Value = '  +5.69%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '4.264.59'
$ws.Cells.Item(15,4).ClearFormats()
$ws.Cells.Item(15,5).Value = '  +4.74%  '

$ws.Cells.Item(16,5).Value = '  -0.84%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '3.814.31'
$ws.Cells.Item(17,4).ClearFormats()
$ws.Cells.Item(17,5).Value = '  +8.46%  '

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '20.47'
$ws.Cells.Item(18,4).ClearFormats()
$ws.Cells.Item(18,5).Value = '  -1.13%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '13.06'
$ws.Cells.Item(19,4).ClearFormats()
$ws.Cells.Item(19,5).Value = '  +4.98%  '

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '1.12'
$ws.Cells.Item(20,4).ClearFormats()
$ws.Cells.Item(20,5).Value = '  +0.68%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '66.689.14'
$ws.Cells.Item(21,4).ClearFormats()
$ws.Cells.Item(21,5).Value = '  +2.70%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '439.55'
$ws.Cells.Item(22,4).ClearFormats()
$ws.Cells.Item(22,5).Value = '  -5.46%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '16.41'
$ws.Cells.Item(23,4).ClearFormats()
$ws.Cells.Item(23,5).Value = '  +20.59%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '89.37'
$ws.Cells.Item(24,4).ClearFormats()
$ws.Cells.Item(24,5).Value = '  -2.70%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '3.13'
$ws.Cells.Item(25,4).ClearFormats()
$ws.Cells.Item(25,5).Value = '  -4.71%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '37.38'
$ws.Cells.Item(26,4).ClearFormats()
$ws.Cells.Item(26,5).Value = '  +8.18%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '10.33'
$ws.Cells.Item(27,4).ClearFormats()
$ws.Cells.Item(27,5).Value = '  +1.11%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '3.30'
$ws.Cells.Item(28,4).ClearFormats()
$ws.Cells.Item(28,5).Value = '  -2.14%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '5.01'
$ws.Cells.Item(29,4).ClearFormats()
$ws.Cells.Item(29,5).Value = '  +4.09%  '

$ws.Cells.Item(30,2).Value = 'Hedera'
$ws.Cells.Item(30,3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '0.126'
$ws.Cells.Item(30,4).ClearFormats()
$ws.Cells.Item(30,5).Value = '  +9.35%  '

$ws.Cells.Item(31,2).Value = 'Cosmos'
$ws.Cells.Item(31,3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '12.76'
$ws.Cells.Item(31,4).ClearFormats()
$ws.Cells.Item(31,5).Value = '  +1.48%  '

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '2.77'
$ws.Cells.Item(32,4).ClearFormats()
$ws.Cells.Item(32,5).Value = '  +2.38%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '7.25'
$ws.Cells.Item(33,4).ClearFormats()
$ws.Cells.Item(33,5).Value = '  -4.87%  '

$ws.Cells.Item(34,5).Value = '  -1.84%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '41.55'
$ws.Cells.Item(35,4).ClearFormats()
$ws.Cells.Item(35,5).Value = '  +1.93%  '

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '57.22'
$ws.Cells.Item(36,4).ClearFormats()
$ws.Cells.Item(36,5).Value = '  -0.99%  '

$ws.Cells.Item(37,5).Value = '  -0.10%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.0493'
$ws.Cells.Item(38,4).ClearFormats()
$ws.Cells.Item(38,5).Value = '  -4.89%  '

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '3.08'
$ws.Cells.Item(39,4).ClearFormats()
$ws.Cells.Item(39,5).Value = '  +31.91%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.0₃0729'
$ws.Cells.Item(40,4).ClearFormats()
$ws.Cells.Item(40,5).Value = '  +1.35%  '

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.149'
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).Value = '  +4.05%  '

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '28.76'
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).Value = '  +29.31%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.996'
$ws.Cells.Item(43,4).ClearFormats()
$ws.Cells.Item(43,5).Value = '  -0.12%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '3.41'
$ws.Cells.Item(44,4).ClearFormats()
$ws.Cells.Item(44,5).Value = '  +0.47%  '

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '148.36'
$ws.Cells.Item(45,4).ClearFormats()
$ws.Cells.Item(45,5).Value = '  +1.39%  '

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '2.11'
$ws.Cells.Item(46,4).ClearFormats()
$ws.Cells.Item(46,5).Value = '  +3.75%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '4.38'
$ws.Cells.Item(47,4).ClearFormats()
$ws.Cells.Item(47,5).Value = '  -3.45%  '

$ws.Cells.Item(48,5).Value = '  -8.49%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '2.59'
$ws.Cells.Item(49,4).ClearFormats()
$ws.Cells.Item(49,5).Value = '  -7.57%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.305'
$ws.Cells.Item(50,4).ClearFormats()
$ws.Cells.Item(50,5).Value = '  -5.34%  '

$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.160'
$ws.Cells.Item(51,4).ClearFormats()
$ws.Cells.Item(51,5).Value = '  +11.22%  '
